$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix pre-existing alignment inconsistencies (rows 10-16) ---
$ws.Range("H10").HorizontalAlignment = -4108

$ws.Range("C3").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("C12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("H11").HorizontalAlignment = -4108

$ws.Range("D14").HorizontalAlignment = -4108
$ws.Range("E14").HorizontalAlignment = -4108
$ws.Range("H14").HorizontalAlignment = -4108

$ws.Range("H15").HorizontalAlignment = -4108
$ws.Range("H16").HorizontalAlignment = -4108

# --- New change-log rows for the final documentation additions ---
$ws.Range("C17").Value = 43474
$ws.Range("C18").Value = 43480
$ws.Range("C19").Value = 43484
$ws.Range("C3").Copy()
$ws.Range("C17:C19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D17").Value = "-"
$ws.Range("E17").Value = "1.0"
$ws.Range("F17").Value = "TODOS"
$ws.Range("G17").Value = "Realización de doc.final"
$ws.Range("H17").Value = "-"
$ws.Range("I17").Value = "Añadido de documentación"

$ws.Range("D18").Value = "1.0"
$ws.Range("E18").Value = "2.0"
$ws.Range("F18").Value = "TODOS"
$ws.Range("G18").Value = "Realización de doc.final"
$ws.Range("H18").Value = "-"
$ws.Range("I18").Value = "Añadido de documentación"

$ws.Range("D19").Value = "2.0"
$ws.Range("E19").Value = "3.0"
$ws.Range("F19").Value = "TODOS"
$ws.Range("G19").Value = "Realización de doc.final"
$ws.Range("H19").Value = "-"
$ws.Range("I19").Value = "Añadido de documentación"

$ws.Range("D17:D19,E17:E19,F17:F19,H17:H19").HorizontalAlignment = -4108

# --- Clear the two now-unused template rows below ---
$ws.Range("A20:B21").ClearContents()

# --- Match the column widths Excel recalculates for the wider text ---
$ws.Columns("F").ColumnWidth = 16.451822916666668
$ws.Columns("G").ColumnWidth = 56.451822916666664

# --- Restore the saved selection ---
$ws.Range("C3:C19").Select()
